$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '43.963.60'
$ws.Range('E2').Value = '  +3.03%  '

# Row 3
$ws.Range('D3').Value = '2.244.46'
$ws.Range('E3').Value = '  +1.92%  '

# Row 4
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '258.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.01%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '80.43'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +8.08%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.626'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.03%  '

# Row 8
$ws.Range('E8').Value = '  -0.05%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.603'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.36%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '43.24'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.25%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0930'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.93%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.08'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.41%  '

# Row 13
$ws.Range('E13').Value = '  +2.07%  '

# Row 14
$ws.Range('D14').Value = '2.575.59'
$ws.Range('E14').Value = '  +1.81%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.71'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.39%  '

# Row 16
$ws.Range('D16').Value = '2.284.55'
$ws.Range('E16').Value = '  +4.09%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.792'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.05%  '

# Row 18
$ws.Range('D18').Value = '43.878.28'
$ws.Range('E18').Value = '  +3.09%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000104'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.71%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.45'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.44%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.05'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.83%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.34'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.83%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.32'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.10%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.35'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.79%  '

# Row 25
$ws.Range('E25').Value = '  +0.10%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.83'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.21%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '40.60'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.74%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.37'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.09%  '

# Row 29
$ws.Range('E29').Value = '  +1.03%  '

# Row 30
$ws.Range('E30').Value = '  -0.60%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '172.68'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.48%  '

# Row 32
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0887'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +10.56%  '

# Row 33
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.61'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.50%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.30'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.02%  '

# Row 35
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.123'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.99%  '

# Row 36
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.112'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.35%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0363'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.40%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.52'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.21%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.84'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.65%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.93'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +20.89%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.15'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.16%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.56'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.53%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '63.04'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.46%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.203'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.71%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '103.92'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.74%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.52'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.71%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0987'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.63%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.450'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -8.26%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.12'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.81%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.15'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.07%  '

# Row 51
$ws.Range('E51').Value = '  +24.87%  '
